$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 13
$ws.Range("B6").Value = 70.83
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 4

$ws.Range("C7").Select()
